$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.272.44'
$ws.Range("E2").Value = '  +1.76%  '
$ws.Range("D3").Value = '2.345.96'
$ws.Range("E3").Value = '  +0.59%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '''542.66'
$ws.Range("E5").Value = '  +1.78%  '
$ws.Range("D6").Value = '''135.02'
$ws.Range("E6").Value = '  +1.62%  '
$ws.Range("E7").Value = '  +0.73%  '
$ws.Range("E8").Value = '  +4.85%  '
$ws.Range("E9").Value = '  +0.26%  '
$ws.Range("E10").Value = '  +6.66%  '
$ws.Range("E11").Value = '  -0.55%  '
$ws.Range("E12").Value = '  +3.37%  '
$ws.Range("E13").Value = '  +0.87%  '
$ws.Range("D14").Value = '2.763.22'
$ws.Range("E14").Value = '  +0.47%  '
$ws.Range("D15").Value = '58.215.73'
$ws.Range("E15").Value = '  +1.64%  '
$ws.Range("E16").Value = '  +0.21%  '
$ws.Range("D17").Value = '2.345.51'
$ws.Range("E17").Value = '  -0.02%  '
$ws.Range("D18").Value = '''10.72'
$ws.Range("E18").Value = '  +2.48%  '
$ws.Range("D19").Value = '''333.38'
$ws.Range("E19").Value = '  -1.87%  '
$ws.Range("E20").Value = '  +1.82%  '
$ws.Range("E21").Value = '  -3.75%  '
$ws.Range("E22").Value = '  +0.08%  '
$ws.Range("D23").Value = '''5.62'
$ws.Range("E23").Value = '  +0.61%  '
$ws.Range("D24").Value = '''62.80'
$ws.Range("E24").Value = '  +1.55%  '
$ws.Range("E25").Value = '  +1.78%  '
$ws.Range("E26").Value = '  -4.23%  '
$ws.Range("D27").Value = '''1.00'
$ws.Range("D28").Value = '''1.41'
$ws.Range("E28").Value = '  +5.35%  '
$ws.Range("E29").Value = '  +1.90%  '
$ws.Range("D30").Value = '''170.38'
$ws.Range("E30").Value = '  +0.31%  '
$ws.Range("E31").Value = '  +1.16%  '
$ws.Range("E32").Value = '  -0.43%  '
$ws.Range("E33").Value = '  +12.70%  '
$ws.Range("E34").Value = '  -0.56%  '
$ws.Range("D36").Value = '''4.26'
$ws.Range("E36").Value = '  +5.68%  '
$ws.Range("D37").Value = '''1.00'
$ws.Range("E37").Value = '  +0.96%  '
$ws.Range("E38").Value = '  -2.04%  '
$ws.Range("E39").Value = '  +3.53%  '
$ws.Range("D40").Value = '''39.14'
$ws.Range("E40").Value = '  +0.43%  '
$ws.Range("D41").Value = '''142.45'
$ws.Range("E41").Value = '  -3.78%  '
$ws.Range("E42").Value = '  +1.43%  '
$ws.Range("E43").Value = '  -0.17%  '
$ws.Range("D44").Value = '''288.08'
$ws.Range("E44").Value = '  +0.33%  '
$ws.Range("D45").Value = '''0.0937'
$ws.Range("E45").Value = '  +0.55%  '
$ws.Range("D46").Value = '''19.18'
$ws.Range("E46").Value = '  +1.95%  '
$ws.Range("D47").Value = '''0.0504'
$ws.Range("E47").Value = '  -0.25%  '
$ws.Range("E48").Value = '  +0.37%  '
$ws.Range("D49").Value = '''0.0219'
$ws.Range("E49").Value = '  +0.51%  '
$ws.Range("D50").Value = '''0.382'
$ws.Range("E50").Value = '  +0.97%  '
$ws.Range("D51").Value = '''17.47'
$ws.Range("E51").Value = '  +0.55%  '
